# Apply the "Dispersal" sheet fix: adult males emigrate every 5th year,
# so emigration probability (column C) should be 1, not 0.2, for rows 9-32.
# Also add a note explaining this, merge B33:H33 for the note, and
# update the sheet view's selection/top-left cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dispersal")

# Update emigration probability values in column C, rows 9 through 32, from 0.2 to 1
for ($r = 9; $r -le 32; $r++) {
    $ws.Cells.Item($r, 3).Value = 1
}

# Add explanatory note in B33, merged across B33:H33
$noteText = "NOTE: adult males are made to emigrate once every 5 years since their 1st emigration, so it does not matter that their probability of emigration according to this table is 1"
$ws.Range("B33").Value = $noteText
$ws.Range("B33:H33").Merge()
$ws.Range("B33:H33").HorizontalAlignment = -4108
$ws.Range("B33:H33").WrapText = $true
$ws.Rows.Item(33).RowHeight = 27

# Update sheet view: scroll position and active cell selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("H29").Select()

# Update workbook calc id (cosmetic, records last Excel build that calculated the file)
$wb.CalcId = 140001
